$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# South America (Sao Paulo, Brazil) surprise songs - finish filling in
# dress color, instrument, and song for rows 128-133.

$ws.Cells.Item(128, 5).Value = "blue"
$ws.Cells.Item(128, 6).Value = "guitar"
$ws.Cells.Item(128, 7).Value = "Now That We Don't Talk (Taylor's Version) [From The Vault]"

$ws.Cells.Item(129, 5).Value = "blue"
$ws.Cells.Item(129, 6).Value = "piano"
$ws.Cells.Item(129, 7).Value = "Innocent (Taylor's Version)"

$ws.Cells.Item(130, 5).Value = "yellow"
$ws.Cells.Item(130, 6).Value = "guitar"
$ws.Cells.Item(130, 7).Value = "Safe & Sound (Taylor's Version)"

$ws.Cells.Item(131, 5).Value = "yellow"
$ws.Cells.Item(131, 6).Value = "piano"
$ws.Cells.Item(131, 7).Value = "Untouchable (Taylor's Version)"

$ws.Cells.Item(132, 5).Value = "green"
$ws.Cells.Item(132, 6).Value = "guitar"
$ws.Cells.Item(132, 7).Value = "Say Don't Go (Taylor's Version) [From The Vault]"

$ws.Cells.Item(133, 5).Value = "green"
$ws.Cells.Item(133, 6).Value = "piano"
$ws.Cells.Item(133, 7).Value = "it's time to go"

# Update the selected cell shown in the sheet view to match the new last row
[void]$ws.Range("G133").Select()
